$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rename (row 1): Spanish labels -> snake_case field names ---
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# --- Spanish preposition title-casing fix (de/del/el/la/los/las/y -> De/Del/El/La/Los/Las/Y) ---
$ws.Range("B4").Value = 'Rincón De Romos'
$ws.Range("B5").Value = 'San Francisco De Los Romo'
$ws.Range("B22").Value = 'Amatenango De La Frontera'
$ws.Range("B30").Value = 'Comitán De Domínguez'
$ws.Range("B46").Value = 'Ocozocoautla De Espinosa'
$ws.Range("B78").Value = 'San Juan De Sabinas'
$ws.Range("A81").Value = 'Ciudad De México'
$ws.Range("B85").Value = 'Cuajimalpa De Morelos'
$ws.Range("B106").Value = 'San Juan Del Río'
$ws.Range("A110").Value = 'Estado De México'
$ws.Range("B110").Value = 'Acambay De Ruíz Castañeda'
$ws.Range("B112").Value = 'Almoloya De Alquisiras'
$ws.Range("B113").Value = 'Almoloya De Juárez'
$ws.Range("B116").Value = 'Atizapán De Zaragoza'
$ws.Range("B121").Value = 'Coacalco De Berriozábal'
$ws.Range("B125").Value = 'Ecatepec De Morelos'
$ws.Range("B132").Value = 'Naucalpan De Juárez'
$ws.Range("B137").Value = 'San Felipe Del Progreso'
$ws.Range("B138").Value = 'Soyaniquilpan De Juárez'
$ws.Range("B150").Value = 'Tlalnepantla De Baz'
$ws.Range("B160").Value = 'San Miguel De Allende'
$ws.Range("B161").Value = 'Apaseo El Alto'
$ws.Range("B168").Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range("B181").Value = 'San Luis De La Paz'
$ws.Range("B185").Value = 'Valle De Santiago'
$ws.Range("B189").Value = 'Acapulco De Juárez'
$ws.Range("B192").Value = 'Ajuchitlán Del Progreso'
$ws.Range("B195").Value = 'Atlamajalcingo Del Monte'
$ws.Range("B197").Value = 'Ayutla De Los Libres'
$ws.Range("B199").Value = 'Chilpancingo De Los Bravo'
$ws.Range("B200").Value = 'Coahuayutla De José María Izazaga'
$ws.Range("B203").Value = 'Coyuca De Benítez'
$ws.Range("B204").Value = 'Coyuca De Catalán'
$ws.Range("B207").Value = 'Cuetzala Del Progreso'
$ws.Range("B208").Value = 'Cutzamala De Pinzón'
$ws.Range("B213").Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range("B215").Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range("B217").Value = 'Mártir De Cuilapan'
$ws.Range("B227").Value = 'Taxco De Alarcón'
$ws.Range("B229").Value = 'Técpan De Galeana'
$ws.Range("B234").Value = 'Tlalixtaquilla De Maldonado'
$ws.Range("B235").Value = 'Tlapa De Comonfort'
$ws.Range("B248").Value = 'Atotonilco El Grande'
$ws.Range("B252").Value = 'Cuautepec De Hinojosa'
$ws.Range("B254").Value = 'Huejutla De Reyes'
$ws.Range("B257").Value = 'Jacala De Ledezma'
$ws.Range("B262").Value = 'Mixquiahuala De Juárez'
$ws.Range("B264").Value = 'Pachuca De Soto'
$ws.Range("B266").Value = 'Progreso De Obregón'
$ws.Range("B272").Value = 'Tenango De Doria'
$ws.Range("B273").Value = 'Tepehuacán De Guerrero'
$ws.Range("B277").Value = 'Tulancingo De Bravo'
$ws.Range("B278").Value = 'Villa De Tezontepec'
$ws.Range("B280").Value = 'Zacualtipán De Ángeles'
$ws.Range("B284").Value = 'Acatlán De Juárez'
$ws.Range("B293").Value = 'Encarnación De Díaz'
$ws.Range("B298").Value = 'Jilotlán De Los Dolores'
$ws.Range("B299").Value = 'Lagos De Moreno'
$ws.Range("B302").Value = 'Ojuelos De Jalisco'
$ws.Range("B305").Value = 'San Juan De Los Lagos'
$ws.Range("B308").Value = 'Tepatitlán De Morelos'
$ws.Range("B310").Value = 'Tizapán El Alto'
$ws.Range("B359").Value = 'Coatlán Del Río'
$ws.Range("B365").Value = 'Jonacatepec De Leandro Valle'
$ws.Range("B369").Value = 'Tetela Del Volcán'
$ws.Range("B373").Value = 'Zacualpan De Amilpas'
$ws.Range("B384").Value = 'San Nicolás De Los Garza'
$ws.Range("B387").Value = 'Acatlán De Pérez Figueroa'
$ws.Range("B391").Value = 'Chalcatongo De Hidalgo'
$ws.Range("B392").Value = 'Coicoyán De Las Flores'
$ws.Range("B395").Value = 'El Barrio De La Soledad'
$ws.Range("B396").Value = 'Fresnillo De Trujano'
$ws.Range("B397").Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range("B398").Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range("B399").Value = 'Ixtlán De Juárez'
$ws.Range("B400").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B403").Value = 'Mazatlán Villa De Flores'
$ws.Range("B404").Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range("B405").Value = 'Oaxaca De Juárez'
$ws.Range("B406").Value = 'Ocotlán De Morelos'
$ws.Range("B407").Value = 'Putla Villa De Guerrero'
$ws.Range("B416").Value = 'San Felipe Jalapa De Díaz'
$ws.Range("B477").Value = 'Tanetze De Zaragoza'
$ws.Range("B478").Value = 'Tataltepec De Valdés'
$ws.Range("B479").Value = 'Teococuilco De Marcos Pérez'
$ws.Range("B480").Value = 'Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca'
$ws.Range("B481").Value = 'Tlacolula De Matamoros'
$ws.Range("B482").Value = 'Totontepec Villa De Morelos'
$ws.Range("B483").Value = 'Villa De Tututepec'
$ws.Range("B484").Value = 'Villa Sola De Vega'
$ws.Range("B506").Value = 'Cuayuca De Andrade'
$ws.Range("B507").Value = 'Cuetzalan Del Progreso'
$ws.Range("B517").Value = 'Huehuetlán El Chico'
$ws.Range("B521").Value = 'Izúcar De Matamoros'
$ws.Range("B525").Value = 'Los Reyes De Juárez'
$ws.Range("B529").Value = 'Palmar De Bravo'
$ws.Range("B541").Value = 'San Salvador El Verde'
$ws.Range("B546").Value = 'Tepanco De López'
$ws.Range("B547").Value = 'Tepatlaxco De Hidalgo'
$ws.Range("B551").Value = 'Tepexi De Rodríguez'
$ws.Range("B555").Value = 'Tlacotepec De Benito Juárez'
$ws.Range("B574").Value = 'Amealco De Bonfil'
$ws.Range("B576").Value = 'Cadereyta De Montes'
$ws.Range("B578").Value = 'Jalpan De Serra'
$ws.Range("B579").Value = 'Landa De Matamoros'
$ws.Range("B581").Value = 'Pinal De Amoles'
$ws.Range("B583").Value = 'San Juan Del Río'
$ws.Range("B589").Value = 'Axtla De Terrazas'
$ws.Range("B592").Value = 'Ciudad Del Maíz'
$ws.Range("B597").Value = 'Mexquitic De Carmona'
$ws.Range("B602").Value = 'Santa María Del Río'
$ws.Range("B607").Value = 'Villa De Guadalupe'
$ws.Range("B608").Value = 'Villa De Ramos'
$ws.Range("B624").Value = 'Jalpa De Méndez'
$ws.Range("B647").Value = 'San Pablo Del Monte'
$ws.Range("B664").Value = 'Amatlán De Los Reyes'
$ws.Range("B668").Value = 'Castillo De Teayo'
$ws.Range("B670").Value = 'Cazones De Herrera'
$ws.Range("B681").Value = 'Cosamaloapan De Carpio'
$ws.Range("B688").Value = 'Hueyapan De Ocampo'
$ws.Range("B689").Value = 'Ignacio De La Llave'
$ws.Range("B691").Value = 'Ixhuatlán De Madero'
$ws.Range("B692").Value = 'Ixhuatlán Del Café'
$ws.Range("B693").Value = 'Ixhuatlán Del Sureste'
$ws.Range("B698").Value = 'Juchique De Ferrer'
$ws.Range("B702").Value = 'Martínez De La Torre'
$ws.Range("B704").Value = 'Medellín De Bravo'
$ws.Range("B707").Value = 'Mixtla De Altamirano'
$ws.Range("B713").Value = 'Paso De Ovejas'
$ws.Range("B719").Value = 'Sayula De Alemán'
$ws.Range("B721").Value = 'Soledad De Doblado'
$ws.Range("B738").Value = 'Vega De Alatorre'
$ws.Range("B754").Value = 'Mezquital Del Oro'
$ws.Range("B756").Value = 'Nochistlán De Mejía'
$ws.Range("B757").Value = 'Noria De Ángeles'

# --- Tiny floating point re-computation fixes on a handful of subtotal rows ---
$ws.Range("D67").Value = 0.09525807805287452
$ws.Range("D242").Value = 0.09441879983214437
$ws.Range("D638").Value = 0.009232060428031891
$ws.Range("D743").Value = 0.0986151909357952

# --- Drop trailing footnote/metadata rows (769-773); shrinks used range to A1:D767 ---
$ws.Rows("769:773").Delete()
